# Update the Ecopulpers "Productivity increase due to 100% new machines"
# sensitivity results table. The underlying sensitivity model was
# re-run in "large scale" intervention mode, which changes the scaled
# results (Investment, Saving, PROI, PPBT and the various
# Water/Emission/Land/Import/Workforce/Capital Saving, Investment and
# Total Impact columns) for every scenario row (rows 4-12, columns D:Y).
# Columns A:C (labels) and rows 1-3 (headers) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = 96.95999200548977
$ws.Range("E4").Value = -2.42207265086472
$ws.Range("F4").Value = -0.02498012428391687
$ws.Range("G4").Value = -40.03182644867132
$ws.Range("H4").Value = 138.7039714279817
$ws.Range("I4").Value = -0.01406058689462952
$ws.Range("J4").Value = 0.05149963283929537
$ws.Range("K4").Value = 0.08904049545526505
$ws.Range("L4").Value = -0.4935806735884398
$ws.Range("M4").Value = -0.4535888531245291
$ws.Range("N4").Value = 0.2179608714068308
$ws.Range("O4").Value = 0.1999594387507386
$ws.Range("P4").Value = 0.0008709411613381235
$ws.Range("Q4").Value = 28.30095673212782
$ws.Range("R4").Value = 5.496810862794518
$ws.Range("S4").Value = 51.23543111188337
$ws.Range("T4").Value = -1386.821753408411
$ws.Range("U4").Value = 0.3405653076970339
$ws.Range("V4").Value = -0.5141253872316156
$ws.Range("W4").Value = 27.41055177757517
$ws.Range("X4").Value = 10.03269939403981
$ws.Range("Y4").Value = 56.17123784776777
$ws.Range("D5").Value = 96.95999200548977
$ws.Range("E5").Value = 8.343645698390901
$ws.Range("F5").Value = 0.08605245860497279
$ws.Range("G5").Value = 11.6208184659841
$ws.Range("H5").Value = 139.1871222900954
$ws.Range("I5").Value = -0.006381419016179279
$ws.Range("J5").Value = 0.1031459945506867
$ws.Range("K5").Value = 1.036461042473093
$ws.Range("L5").Value = -0.1079261575359851
$ws.Range("M5").Value = 0.7617003140039742
$ws.Range("N5").Value = 0.2179608714068308
$ws.Range("O5").Value = 0.1999594387507386
$ws.Range("P5").Value = 0.0008709411613381235
$ws.Range("Q5").Value = 28.30095673212782
$ws.Range("R5").Value = 5.496810862794518
$ws.Range("S5").Value = 51.23543111188337
$ws.Range("T5").Value = -1391.653262029547
$ws.Range("U5").Value = 0.2637736289125314
$ws.Range("V5").Value = -1.030589004345529
$ws.Range("W5").Value = 17.93634630739689
$ws.Range("X5").Value = -2.120192277245224
$ws.Range("Y5").Value = 52.31469268724322
$ws.Range("D6").Value = 96.95999200548977
$ws.Range("E6").Value = 19.10913966968656
$ws.Range("F6").Value = 0.1970827273645466
$ws.Range("G6").Value = 5.074011372647013
$ws.Range("H6").Value = 139.6702630824293
$ws.Range("I6").Value = 0.001297588816669304
$ws.Range("J6").Value = 0.1547912798532707
$ws.Range("K6").Value = 1.983861843356863
$ws.Range("L6").Value = 0.2777203212026507
$ws.Range("M6").Value = 1.976964151952416
$ws.Range("N6").Value = 0.2179608714068308
$ws.Range("O6").Value = 0.1999594387507386
$ws.Range("P6").Value = 0.0008709411613381235
$ws.Range("Q6").Value = 28.30095673212782
$ws.Range("R6").Value = 5.496810862794518
$ws.Range("S6").Value = 51.23543111188337
$ws.Range("T6").Value = -1396.484669952886
$ws.Range("U6").Value = 0.1869835505840456
$ws.Range("V6").Value = -1.547041857371369
$ws.Range("W6").Value = 8.462338298559189
$ws.Range("X6").Value = -14.27283065672964
$ws.Range("Y6").Value = 48.45822789985687
$ws.Range("D7").Value = 96.95999200548977
$ws.Range("E7").Value = 29.87440926767886
$ws.Range("F7").Value = 0.3081106820428307
$ws.Range("G7").Value = 3.245586921458924
$ws.Range("H7").Value = 140.1533938052598
$ws.Range("I7").Value = 0.008976436611192185
$ws.Range("J7").Value = 0.2064354887788795
$ws.Range("K7").Value = 2.931242898805067
$ws.Range("L7").Value = 0.6633587619289756
$ws.Range("M7").Value = 3.192202662117779
$ws.Range("N7").Value = 0.2179608714068308
$ws.Range("O7").Value = 0.1999594387507386
$ws.Range("P7").Value = 0.0008709411613381235
$ws.Range("Q7").Value = 28.30095673212782
$ws.Range("R7").Value = 5.496810862794518
$ws.Range("S7").Value = 51.23543111188337
$ws.Range("T7").Value = -1401.315977181192
$ws.Range("U7").Value = 0.1101950726388168
$ws.Range("V7").Value = -2.063483946627457
$ws.Range("W7").Value = -1.011472255922854
$ws.Range("X7").Value = -26.42521575838327
$ws.Range("Y7").Value = 44.60184349259362
$ws.Range("D8").Value = 96.95999200548977
$ws.Range("E8").Value = 40.639454504475
$ws.Range("F8").Value = 0.4191363227646929
$ws.Range("G8").Value = 2.385858599426157
$ws.Range("H8").Value = 140.6365144589217
$ws.Range("I8").Value = 0.01665512435647543
$ws.Range("J8").Value = 0.2580786213629835
$ws.Range("K8").Value = 3.878604209516197
$ws.Range("L8").Value = 1.048989166272804
$ws.Range("M8").Value = 4.407415844500065
$ws.Range("N8").Value = 0.2179608714068308
$ws.Range("O8").Value = 0.1999594387507386
$ws.Range("P8").Value = 0.0008709411613381235
$ws.Range("Q8").Value = 28.30095673212782
$ws.Range("R8").Value = 5.496810862794518
$ws.Range("S8").Value = 51.23543111188337
$ws.Range("T8").Value = -1406.147183717811
$ws.Range("U8").Value = 0.03340819518598437
$ws.Range("V8").Value = -2.579915272468497
$ws.Range("W8").Value = -10.48508536303416
$ws.Range("X8").Value = -38.57734758220613
$ws.Range("Y8").Value = 40.74553944915533
$ws.Range("D9").Value = 96.95999200548977
$ws.Range("E9").Value = 51.40427538286895
$ws.Range("F9").Value = 0.5301596495589489
$ws.Range("G9").Value = 1.886224273823784
$ws.Range("H9").Value = 141.1196250437642
$ws.Range("I9").Value = 0.02433365207616589
$ws.Range("J9").Value = 0.3097206776392341
$ws.Range("K9").Value = 4.825945775955915
$ws.Range("L9").Value = 1.434611532604322
$ws.Range("M9").Value = 5.622603700961918
$ws.Range("N9").Value = 0.2179608714068308
$ws.Range("O9").Value = 0.1999594387507386
$ws.Range("P9").Value = 0.0008709411613381235
$ws.Range("Q9").Value = 28.30095673212782
$ws.Range("R9").Value = 5.496810862794518
$ws.Range("S9").Value = 51.23543111188337
$ws.Range("T9").Value = -1410.978289566236
$ws.Range("U9").Value = -0.04337708201092028
$ws.Range("V9").Value = -3.096335835231002
$ws.Range("W9").Value = -19.95850102743134
$ws.Range("X9").Value = -50.72922614682466
$ws.Range("Y9").Value = 36.88931578584015
$ws.Range("D10").Value = 96.95999200548977
$ws.Range("E10").Value = 62.16887191031128
$ws.Range("F10").Value = 0.6411806625024407
$ws.Range("G10").Value = 1.559622830946174
$ws.Range("H10").Value = 141.6027255600638
$ws.Range("I10").Value = 0.03201201977026358
$ws.Range("J10").Value = 0.3613616576394634
$ws.Range("K10").Value = 5.773267598822713
$ws.Range("L10").Value = 1.820225863019004
$ws.Range("M10").Value = 6.837766231037676
$ws.Range("N10").Value = 0.2179608714068308
$ws.Range("O10").Value = 0.1999594387507386
$ws.Range("P10").Value = 0.0008709411613381235
$ws.Range("Q10").Value = 28.30095673212782
$ws.Range("R10").Value = 5.496810862794518
$ws.Range("S10").Value = 51.23543111188337
$ws.Range("T10").Value = -1415.809294729232
$ws.Range("U10").Value = -0.1201607589518972
$ws.Range("V10").Value = -3.612745635233296
$ws.Range("W10").Value = -29.43171925609931
$ws.Range("X10").Value = -62.88085144758224
$ws.Range("Y10").Value = 33.03317248169333
$ws.Range("D11").Value = 96.95999200548977
$ws.Range("E11").Value = 72.93324409518391
$ws.Range("F11").Value = 0.7521993616816152
$ws.Range("G11").Value = 1.329434789421254
$ws.Range("H11").Value = 142.0858160081116
$ws.Range("I11").Value = 0.03969022742603556
$ws.Range("J11").Value = 0.4130015614000513
$ws.Range("K11").Value = 6.72056967834942
$ws.Range("L11").Value = 2.205832156352699
$ws.Range("M11").Value = 8.052903436124325
$ws.Range("N11").Value = 0.2179608714068308
$ws.Range("O11").Value = 0.1999594387507386
$ws.Range("P11").Value = 0.0008709411613381235
$ws.Range("Q11").Value = 28.30095673212782
$ws.Range("R11").Value = 5.496810862794518
$ws.Range("S11").Value = 51.23543111188337
$ws.Range("T11").Value = -1420.640199209709
$ws.Range("U11").Value = -0.196942835509617
$ws.Range("V11").Value = -4.129144672839175
$ws.Range("W11").Value = -38.90474005136639
$ws.Range("X11").Value = -75.03222349844873
$ws.Range("Y11").Value = 29.17710954835638
$ws.Range("D12").Value = 96.95999200548977
$ws.Range("E12").Value = -13.18801538553089
$ws.Range("F12").Value = -0.1360150213789642
$ws.Range("G12").Value = -7.352129124134064
$ws.Range("H12").Value = 138.2208104957681
$ws.Range("I12").Value = -0.02173991483141435
$ws.Range("J12").Value = -0.0001478053145547165
$ws.Range("K12").Value = -0.8583997981622815
$ws.Range("L12").Value = -0.8792432283516973
$ws.Range("M12").Value = -1.668903349898756
$ws.Range("N12").Value = 0.2179608714068308
$ws.Range("O12").Value = 0.1999594387507386
$ws.Range("P12").Value = 0.0008709411613381235
$ws.Range("Q12").Value = 28.30095673212782
$ws.Range("R12").Value = 5.496810862794518
$ws.Range("S12").Value = 51.23543111188337
$ws.Range("T12").Value = -1381.990144086274
$ws.Range("U12").Value = 0.4173585870648822
$ws.Range("V12").Value = 0.002348994306885288
$ws.Range("W12").Value = 36.88495471375063
$ws.Range("X12").Value = 22.18584436178207
$ws.Range("Y12").Value = 60.02786339540035
